$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record at row 6 (pushes existing rows 6..49 down
# to 7..50 -- the former last row, 49, lands on row 50 with no further
# changes needed).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44547
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112031
$ws.Range("G6").Value = "Poroto verde"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 580
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
